# Refresh cryptos list (price + 1h volume %) as produced by the
# scheduled GitHub Actions scraper run. D-column values that look like
# plain numbers are entered with a leading apostrophe so Excel keeps
# them as text (matching the original inlineStr cells) instead of
# silently converting them to numeric cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.466.16"
$ws.Range("E2").Value = "  -2.46%  "
$ws.Range("D3").Value = "3.181.83"
$ws.Range("E3").Value = "  -4.04%  "
$ws.Range("D5").Value = "'572.94"
$ws.Range("E5").Value = "  -2.22%  "
$ws.Range("D6").Value = "'169.36"
$ws.Range("E6").Value = "  -6.90%  "
$ws.Range("D7").Value = "'0.609"
$ws.Range("E7").Value = "  -6.83%  "
$ws.Range("E8").Value = "  -0.17%  "
$ws.Range("D9").Value = "3.190.41"
$ws.Range("E9").Value = "  -3.73%  "
$ws.Range("E10").Value = "  -3.92%  "
$ws.Range("D11").Value = "'6.82"
$ws.Range("E11").Value = "  +0.11%  "
$ws.Range("D12").Value = "'0.390"
$ws.Range("E12").Value = "  -2.68%  "
$ws.Range("D13").Value = "3.746.22"
$ws.Range("E13").Value = "  -3.73%  "
$ws.Range("E14").Value = "  -1.47%  "
$ws.Range("D15").Value = "64.509.93"
$ws.Range("E15").Value = "  -2.46%  "
$ws.Range("D16").Value = "'25.31"
$ws.Range("E16").Value = "  -3.20%  "
$ws.Range("E17").Value = "  -3.67%  "
$ws.Range("D18").Value = "3.178.63"
$ws.Range("E18").Value = "  -2.72%  "
$ws.Range("D19").Value = "'420.26"
$ws.Range("E19").Value = "  -1.15%  "
$ws.Range("D20").Value = "'13.01"
$ws.Range("E20").Value = "  -0.70%  "
$ws.Range("D21").Value = "'5.37"
$ws.Range("E21").Value = "  -3.10%  "
$ws.Range("D22").Value = "'7.17"
$ws.Range("E22").Value = "  -2.72%  "
$ws.Range("E23").Value = "  -0.02%  "
$ws.Range("B24").Value = "Litecoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D24").Value = "'70.31"
$ws.Range("E24").Value = "  -1.89%  "
$ws.Range("B25").Value = "LEO"
$ws.Range("C25").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D25").Value = "'5.67"
$ws.Range("E25").Value = "  -0.03%  "
$ws.Range("E26").Value = "  +2.55%  "
$ws.Range("D27").Value = "'0.497"
$ws.Range("E27").Value = "  -2.80%  "
$ws.Range("E28").Value = "  -7.48%  "
$ws.Range("D29").Value = "'8.75"
$ws.Range("E29").Value = "  -1.38%  "
$ws.Range("D30").Value = "'0.997"
$ws.Range("E30").Value = "  -0.40%  "
$ws.Range("E31").Value = "  -4.41%  "
$ws.Range("D32").Value = "'21.82"
$ws.Range("E32").Value = "  -2.44%  "
$ws.Range("E33").Value = "  -0.11%  "
$ws.Range("D34").Value = "'5.05"
$ws.Range("E34").Value = "  -2.37%  "
$ws.Range("D35").Value = "'6.36"
$ws.Range("E35").Value = "  -2.90%  "
$ws.Range("D36").Value = "'157.24"
$ws.Range("E36").Value = "  -2.27%  "
$ws.Range("E37").Value = "  -4.19%  "
$ws.Range("E38").Value = "  -4.56%  "
$ws.Range("E39").Value = "  -4.96%  "
$ws.Range("D40").Value = "2.698.62"
$ws.Range("E40").Value = "  -5.96%  "
$ws.Range("E41").Value = "  -1.54%  "
$ws.Range("D42").Value = "'24.31"
$ws.Range("E42").Value = "  -7.71%  "
$ws.Range("D43").Value = "'39.30"
$ws.Range("E43").Value = "  -1.28%  "
$ws.Range("D44").Value = "'0.717"
$ws.Range("E44").Value = "  -5.50%  "
$ws.Range("D45").Value = "'0.0621"
$ws.Range("E45").Value = "  -5.71%  "
$ws.Range("D46").Value = "'5.53"
$ws.Range("E46").Value = "  -6.32%  "
$ws.Range("E47").Value = "  -2.76%  "
$ws.Range("D48").Value = "'291.83"
$ws.Range("E48").Value = "  -6.59%  "
$ws.Range("D49").Value = "'21.40"
$ws.Range("E49").Value = "  -7.34%  "
$ws.Range("D50").Value = "'0.0997"
$ws.Range("E50").Value = "  -5.78%  "
$ws.Range("D51").Value = "'0.997"
$ws.Range("E51").Value = "  -0.21%  "
